$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.597.46"
$ws.Range("E2").Value = "  +1.63%  "
$ws.Range("D3").Value = "3.458.56"
$ws.Range("E3").Value = "  +1.87%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.39"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.62"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.458.67"
$ws.Range("E8").Value = "  +1.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.579"
$ws.Range("E9").Value = "  +8.71%  "
$ws.Range("E10").Value = "  -2.16%  "
$ws.Range("E11").Value = "  +4.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.441"
$ws.Range("E12").Value = "  +1.18%  "
$ws.Range("D13").Value = "4.060.63"
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("E14").Value = "  -2.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000194"
$ws.Range("E15").Value = "  +5.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.87"
$ws.Range("E16").Value = "  +6.68%  "
$ws.Range("D17").Value = "64.685.59"
$ws.Range("E17").Value = "  +1.67%  "
$ws.Range("D18").Value = "3.455.15"
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.36"
$ws.Range("E20").Value = "  +2.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "391.39"
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.16"
$ws.Range("E22").Value = "  -3.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.548"
$ws.Range("E23").Value = "  +2.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.93"
$ws.Range("E24").Value = "  +3.01%  "
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000123"
$ws.Range("E26").Value = "  +19.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.45"
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.19"
$ws.Range("E30").Value = "  +9.80%  "
$ws.Range("E31").Value = "  +6.43%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("B33").Value = "RenderToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.58"
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.64"
$ws.Range("E34").Value = "  +2.22%  "
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.07"
$ws.Range("E36").Value = "  +5.39%  "
$ws.Range("E37").Value = "  +1.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "161.91"
$ws.Range("E38").Value = "  +2.15%  "
$ws.Range("E39").Value = "  +1.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0773"
$ws.Range("E40").Value = "  +1.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.54"
$ws.Range("E41").Value = "  +0.34%  "
$ws.Range("D42").Value = "2.931.89"
$ws.Range("E42").Value = "  +1.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.56"
$ws.Range("E43").Value = "  +5.95%  "
$ws.Range("E44").Value = "  -1.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.83"
$ws.Range("E45").Value = "  +3.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.775"
$ws.Range("E46").Value = "  +1.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.23"
$ws.Range("E47").Value = "  +8.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.09"
$ws.Range("E48").Value = "  +2.66%  "
$ws.Range("B49").Value = "SuiNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.873"
$ws.Range("E49").Value = "  +7.04%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.19"
$ws.Range("E50").Value = "  +13.13%  "
$ws.Range("E51").Value = "  +3.83%  "
